$p = $ppt.ActivePresentation

# --- Slide 11: "2 to 3 times" -> "2 times" ---------------------------------
$s11 = $p.Slides.Item(11)
$sh11 = $s11.Shapes.Item(3)
$tr11 = $sh11.TextFrame.TextRange
$para11 = $tr11.Paragraphs(2)
$run11 = $para11.Runs(1)
$run11.Text = "We can also see that the Voodooists and the Others are 2 times more victims than the Catholics and the Protestants."

# --- Slide 15: replace second textbox body + shrink its height -------------
$s15 = $p.Slides.Item(15)
$sh15 = $s15.Shapes.Item(3)
$tr15 = $sh15.TextFrame.TextRange
$tr15.Text = "Among the women who are victims of violence, the women who do not work are more beaten when they go out without telling their spouse than those who work.."
$tr15.Font.Size = 18
$sh15.Height = 159.9468503937008

# --- Slide 5: cosmetic re-touch of the "For every 1,000 women..." run ------
$s5 = $p.Slides.Item(5)
$sh5 = $s5.Shapes.Item(2)
$tr5 = $sh5.TextFrame.TextRange
$para5 = $tr5.Paragraphs(1)
$run5 = $para5.Runs(1)
$run5.Text = $run5.Text

# --- Slide 6: ".." -> "." at the end of the last bullet paragraph ----------
$s6 = $p.Slides.Item(6)
$sh6 = $s6.Shapes.Item(1)
$tr6 = $sh6.TextFrame.TextRange
$para6 = $tr6.Paragraphs(5)
$run6 = $para6.Runs(2)
$run6.Text = "."
